# AShot without inspirepak and jenkins propfile
#
# The ink colour rows for the two job elements (1) 591345 2p and
# 2) 591346 2p) get re-ordered, and the two "Pallet (packaging)" rows
# swap which job element / quantity they belong to.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ink / Varnish rows for "2) 591346  2p" (rows 4-7): Magenta/Black swap
$ws.Range("D4").Value = "Black - Sheet-fed Offset UV - "
$ws.Range("D5").Value = "Magenta - Sheet-fed Offset UV - "

# --- Ink / Varnish rows for "1) 591345  2p" (rows 8-11): rotate Black/Yellow/Cyan
$ws.Range("D8").Value = "Yellow - Sheet-fed Offset UV - "
$ws.Range("D9").Value = "Cyan - Sheet-fed Offset UV - "
$ws.Range("D10").Value = "Black - Sheet-fed Offset UV - "

# --- Pallet (packaging) rows 12-13: swap job element + quantity
$ws.Range("B12").Value = "1) 591345 2p Packed"
$ws.Range("B13").Value = "2) 591346 2p Packed"

# Quantity column holds text that looks numeric ("2.00" / "1.00"); force
# text entry (leading apostrophe) so it round-trips as a string, not a
# literal number.
$ws.Range("E12").Value = "'2.00"
$ws.Range("E13").Value = "'1.00"
